$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-version"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")
# Clear the Constraint(s) value for the "Extension" row (row 2, column AI)
$wsElem.Range("AI2").Value = ""
# The "Fixed Value" for Extension.url (row 5) is the extension's own URL -
# it shares the same text as the Metadata URL, so it must be updated too.
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-version"
